# "Fixed some files for February"
# Two rows (4 and 6) had placeholder/missing geocoded coordinates (0, 0),
# shown with a red-fill "needs fixing" highlight style. This fills in the
# correct latitude/longitude values and restores the normal (bordered,
# un-highlighted) coordinate cell style used by the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the normal coordinate-cell format (bordered, no fill) from F2/G2
# and apply it (format only) to the previously-highlighted cells.
$ws.Range("F2:G2").Copy()
$ws.Range("F4:G4").PasteSpecial(-4122)
$ws.Range("F6:G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the corrected coordinates.
$ws.Range("F4").Value = 37.883069999999996
$ws.Range("G4").Value = -122.302504
$ws.Range("F6").Value = 37.871011000000003
$ws.Range("G6").Value = -122.253404

# Update the active selection left after the fixes.
$ws.Range("E4:G4").Select() | Out-Null

# Set the page to portrait orientation.
$ws.PageSetup.Orientation = 1 | Out-Null
